# Apply the "Player Info" sheet addition + MATCH_CARD_LINK -> MATCH_CODE rework.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new "Player Info" sheet in front of "ODI Batting".
#    NOTE: inserting a sheet shifts every other sheet's position, and sheet
#    handles obtained *before* the insert resolve by position - so grab the
#    "ODI Batting"/"ODI Bowling" handles again (by name) after this step.
# ---------------------------------------------------------------------------
$battingSheetBeforeAdd = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetBeforeAdd)
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centered-top header look used on the other sheets
# by cloning the formatting straight from an existing header cell.
$battingSheetForFormat = $wb.Worksheets.Item("ODI Batting")
$battingSheetForFormat.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

# Data row - ID stays textual (matches "inlineStr" in the source workbook)
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3976"
$playerInfo.Range("B2").Value = "Veerasammy Permaul"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

$playerInfo.Range("A1").Select()

# Match the page margins used throughout the rest of the workbook
$playerInfo.PageSetup.LeftMargin = 54
$playerInfo.PageSetup.RightMargin = 54
$playerInfo.PageSetup.TopMargin = 72
$playerInfo.PageSetup.BottomMargin = 72
$playerInfo.PageSetup.HeaderMargin = 36
$playerInfo.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2) "ODI Batting": MATCH_CARD_LINK column (D) -> MATCH_CODE w/ bare match id.
#    (rows 2..8, in order, map to these match codes)
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$matchCodesInRowOrder = @("3452", "3453", "3454", "3481", "3581", "3583", "4019")
for ($i = 0; $i -lt $matchCodesInRowOrder.Length; $i++) {
    $row = $i + 2
    $cell = $battingSheet.Range("D" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodesInRowOrder[$i]
}

# ---------------------------------------------------------------------------
# 3) "ODI Bowling": MATCH_CARD_LINK column (B) -> MATCH_CODE w/ bare match id.
#    (same match codes, same row order, as "ODI Batting" above)
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

for ($i = 0; $i -lt $matchCodesInRowOrder.Length; $i++) {
    $row = $i + 2
    $cell = $bowlingSheet.Range("B" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodesInRowOrder[$i]
}
